$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "user does not exists"

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "totp not valid"

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "login valid, awaiting totp auth"

$ws.Range("D18:E18").Select()
